$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Mortgage Real Estate Investment Trust...(16)'
$ws.Cells.Item(2, 2).Value = 0.5660140110979992
$ws.Cells.Item(3, 1).Value = 'Energy Equipment & Services(38)'
$ws.Cells.Item(3, 2).Value = 0.456318762427332
$ws.Cells.Item(4, 1).Value = 'Containers & Packaging(12)'
$ws.Cells.Item(4, 2).Value = 0.4271537997120528
$ws.Cells.Item(5, 1).Value = 'Electric Utilities(28)'
$ws.Cells.Item(5, 2).Value = 0.4204824388844848
$ws.Cells.Item(6, 1).Value = 'Marine(15)'
$ws.Cells.Item(6, 2).Value = 0.4187471659623295
$ws.Cells.Item(7, 1).Value = 'Road & Rail(22)'
$ws.Cells.Item(7, 2).Value = 0.4163434739926311
$ws.Cells.Item(8, 1).Value = 'Equity Real Estate Investment Trusts ...(98)'
$ws.Cells.Item(8, 2).Value = 0.4129983619278561
$ws.Cells.Item(9, 1).Value = 'Multi-Utilities(18)'
$ws.Cells.Item(9, 2).Value = 0.4101314264982653
$ws.Cells.Item(10, 1).Value = 'Auto Components(21)'
$ws.Cells.Item(10, 2).Value = 0.4048105841666986
$ws.Cells.Item(11, 1).Value = 'Machinery(86)'
$ws.Cells.Item(11, 2).Value = 0.3956950225213282
$ws.Cells.Item(12, 1).Value = 'Chemicals(52)'
$ws.Cells.Item(12, 2).Value = 0.3856895279690697
$ws.Cells.Item(13, 1).Value = 'Air Freight & Logistics(11)'
$ws.Cells.Item(13, 2).Value = 0.3843055748893925
$ws.Cells.Item(14, 1).Value = 'ETF(303)'
$ws.Cells.Item(14, 2).Value = 0.3702495491863249
$ws.Cells.Item(15, 1).Value = 'Insurance(75)'
$ws.Cells.Item(15, 2).Value = 0.3663951961680489
$ws.Cells.Item(16, 1).Value = 'Consumer Finance(15)'
$ws.Cells.Item(16, 2).Value = 0.3459232203685237
$ws.Cells.Item(17, 1).Value = 'Oil, Gas & Consumable Fuels(125)'
$ws.Cells.Item(17, 2).Value = 0.3409126068671584
$ws.Cells.Item(18, 1).Value = 'Life Sciences Tools & Services(19)'
$ws.Cells.Item(18, 2).Value = 0.3230085395246081
$ws.Cells.Item(19, 1).Value = 'Building Products(24)'
$ws.Cells.Item(19, 2).Value = 0.3099121279876063
$ws.Cells.Item(20, 1).Value = 'Metals & Mining(106)'
$ws.Cells.Item(20, 2).Value = 0.2832846768859565
$ws.Cells.Item(21, 1).Value = 'Trading Companies & Distributors(25)'
$ws.Cells.Item(21, 2).Value = 0.2783591086449008
$ws.Cells.Item(22, 1).Value = 'Electrical Equipment(28)'
$ws.Cells.Item(22, 2).Value = 0.2747318155666205
$ws.Cells.Item(23, 1).Value = 'Banks(251)'
$ws.Cells.Item(23, 2).Value = 0.2743789436570588
$ws.Cells.Item(24, 1).Value = 'Diversified Telecommunication Services(20)'
$ws.Cells.Item(24, 2).Value = 0.2696989538733799
$ws.Cells.Item(25, 1).Value = 'Capital Markets(76)'
$ws.Cells.Item(25, 2).Value = 0.2692032346221656
$ws.Cells.Item(26, 1).Value = 'Media(42)'
$ws.Cells.Item(26, 2).Value = 0.244822527871412
$ws.Cells.Item(27, 1).Value = 'Aerospace & Defense(37)'
$ws.Cells.Item(27, 2).Value = 0.2389711178610557
$ws.Cells.Item(28, 1).Value = 'Semiconductors & Semiconductor Equipment(70)'
$ws.Cells.Item(28, 2).Value = 0.2340788241204227
$ws.Cells.Item(29, 1).Value = 'Construction & Engineering(21)'
$ws.Cells.Item(29, 2).Value = 0.2257018002790304
$ws.Cells.Item(30, 1).Value = 'Commercial Services & Supplies(52)'
$ws.Cells.Item(30, 2).Value = 0.2169784159426315
$ws.Cells.Item(31, 1).Value = 'Biotechnology(128)'
$ws.Cells.Item(31, 2).Value = 0.2003798688164793
$ws.Cells.Item(32, 1).Value = 'Household Durables(39)'
$ws.Cells.Item(32, 2).Value = 0.1903255716778237
$ws.Cells.Item(33, 1).Value = 'IT Services(52)'
$ws.Cells.Item(33, 2).Value = 0.1894672843098727
$ws.Cells.Item(34, 1).Value = 'Hotels, Restaurants & Leisure(51)'
$ws.Cells.Item(34, 2).Value = 0.1810712110532062
$ws.Cells.Item(35, 1).Value = 'Health Care Providers & Services(47)'
$ws.Cells.Item(35, 2).Value = 0.1787980920664778
$ws.Cells.Item(36, 1).Value = 'Thrifts & Mortgage Finance(47)'
$ws.Cells.Item(36, 2).Value = 0.1757322938026915
$ws.Cells.Item(37, 1).Value = 'Pharmaceuticals(53)'
$ws.Cells.Item(37, 2).Value = 0.1756744816524684
$ws.Cells.Item(38, 1).Value = 'Health Care Equipment & Supplies(86)'
$ws.Cells.Item(38, 2).Value = 0.1687022434171843
$ws.Cells.Item(39, 1).Value = 'Software(70)'
$ws.Cells.Item(39, 2).Value = 0.1674105163098936
$ws.Cells.Item(40, 1).Value = 'Electronic Equipment, Instruments & C...(78)'
$ws.Cells.Item(40, 2).Value = 0.1290376834206546
$ws.Cells.Item(41, 1).Value = 'Specialty Retail(59)'
$ws.Cells.Item(41, 2).Value = 0.1251773521755793
$ws.Cells.Item(42, 1).Value = 'Professional Services(35)'
$ws.Cells.Item(42, 2).Value = 0.1221762510538849
$ws.Cells.Item(43, 1).Value = 'Communications Equipment(45)'
$ws.Cells.Item(43, 2).Value = 0.1031403075966969

Write-Output $ws.UsedRange.Address()
